$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44 / 45: coin identity swap (EnergySwap <-> Decentraland) ---
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6456"
$ws.Range("E44").Value = "  +6.68%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.14"
$ws.Range("E45").Value = "  +5.10%  "

# --- Price (D) / Volume(1h) (E) updates for remaining rows ---
$ws.Range("D2").Value = "30.848.59"
$ws.Range("E2").Value = "  +2.32%  "
$ws.Range("D3").Value = "2.120.54"
$ws.Range("E3").Value = "  +10.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.55"
$ws.Range("E5").Value = "  +4.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5401"
$ws.Range("E7").Value = "  +6.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4431"
$ws.Range("E8").Value = "  +8.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09100"
$ws.Range("E9").Value = "  +9.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.46"
$ws.Range("E10").Value = "  +9.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.187"
$ws.Range("E11").Value = "  +6.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.32"
$ws.Range("E12").Value = "  +5.37%  "
$ws.Range("D13").Value = "2.123.56"
$ws.Range("E13").Value = "  +10.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.813"
$ws.Range("E14").Value = "  +5.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.866"
$ws.Range("E15").Value = "  +8.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "98.48"
$ws.Range("E16").Value = "  +6.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001143"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06659"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.30"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.423"
$ws.Range("E21").Value = "  +7.79%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "30.963.11"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.20"
$ws.Range("E24").Value = "  +7.30%  "
$ws.Range("D25").Value = "2.372.44"
$ws.Range("E25").Value = "  +10.70%  "
$ws.Range("E26").Value = "  +3.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.93"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.599"
$ws.Range("E28").Value = "  +14.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.92"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.21"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.175"
$ws.Range("E31").Value = "  +3.60%  "
$ws.Range("E32").Value = "  +3.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.312"
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.009"
$ws.Range("E34").Value = "  +5.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.555"
$ws.Range("E35").Value = "  +28.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02618"
$ws.Range("E36").Value = "  +6.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.628"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.682"
$ws.Range("E38").Value = "  +12.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06777"
$ws.Range("E39").Value = "  +5.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.87"
$ws.Range("E40").Value = "  +12.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2286"
$ws.Range("E41").Value = "  +6.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6870"
$ws.Range("E42").Value = "  +5.56%  "
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.263"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.688"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.293"
$ws.Range("E49").Value = "  +6.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.64"
$ws.Range("E50").Value = "  +7.48%  "
$ws.Range("E51").Value = "  +3.83%  "
